# Update "想去人数" (F column) figures across the "展览", "演出" and
# "全部类型" sheets, as published in the new data snapshot (456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 5901
$ws1.Range("F4").Value  = 1089
$ws1.Range("F5").Value  = 1051
$ws1.Range("F10").Value = 61
$ws1.Range("F13").Value = 2042
$ws1.Range("F14").Value = 1523
$ws1.Range("F15").Value = 1120
$ws1.Range("F18").Value = 429
$ws1.Range("F20").Value = 234
$ws1.Range("F21").Value = 1072
$ws1.Range("F24").Value = 3693
$ws1.Range("F28").Value = 169
$ws1.Range("F30").Value = 524
$ws1.Range("F35").Value = 328
$ws1.Range("F40").Value = 90

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 754
$ws2.Range("F6").Value = 415

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5901
$ws4.Range("F5").Value  = 1089
$ws4.Range("F7").Value  = 754
$ws4.Range("F8").Value  = 1051
$ws4.Range("F11").Value = 415
$ws4.Range("F15").Value = 61
$ws4.Range("F19").Value = 2042
$ws4.Range("F20").Value = 1523
$ws4.Range("F21").Value = 1120
$ws4.Range("F24").Value = 429
$ws4.Range("F27").Value = 234
$ws4.Range("F28").Value = 1072
$ws4.Range("F30").Value = 3693
$ws4.Range("F34").Value = 169
$ws4.Range("F36").Value = 524
$ws4.Range("F41").Value = 328
$ws4.Range("F46").Value = 90

$wb.Save()
